# Update countries & provincias Spain
# Refresh the "Casos totales"-and-friends stats for a handful of countries,
# then re-sort the whole ranking table (A4:H219) by "Casos totales"
# (column B) descending - exactly like the live dashboard export would -
# and bump the "updated at" footer timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# country name -> updated column values (only columns that actually change)
$updates = @{
    "Estados Unidos"       = @{ B = 2462713; C = 159;  E = 1297823 }
    "Banglades"            = @{ B = 126606;  C = 3946; D = 51495; E = 73490; G = 39; H = 1621 }
    "Belgica"              = @{ B = 61007;   C = 109;  E = 34510; G = 4;  H = 9726 }
    "Indonesia"            = @{ B = 50187;   C = 1178; D = 20449; E = 27118; G = 47; H = 2620 }
    "Oman"                 = @{ B = 34902;   C = 1366; D = 18520; E = 16238; G = 2;  H = 144 }
    "Polonia"              = @{ B = 33119;   C = 298;  E = 13053; G = 16; H = 1412 }
    "Filipinas"            = @{ B = 33069;   C = 774;  D = 8910;  E = 22947; G = 8;  H = 1212 }
    "Israel"               = @{ B = 22139;   C = 95;   D = 15961; E = 5870 }
    "Austria"              = @{ B = 17477;   C = 28;   D = 16320; E = 459;  G = 5;  H = 698 }
    "Malasia"              = @{ B = 8600;    C = 4;    D = 8271;  E = 208 }
    "El Salvador"          = @{ B = 5336;    C = 186;  D = 3116;  E = 2094; G = 7;  H = 126 }
    "Albania"              = @{ B = 2192;    C = 78;   D = 1250;  E = 894;  G = 1;  H = 48 }
    "Eslovenia"            = @{ B = 1547;    C = 6;    E = 62 }
    "Estado de Palestina"  = @{ B = 1362;    C = 34;   E = 917 }
}

$firstRow = 4
$lastRow = 219

# Map column letter -> column index without relying on Range.Find (which
# does a *substring* search by default, e.g. "Estados Unidos" would also
# hit "Islas Virgenes de los Estados Unidos") - walk column A and match
# whole cell text instead, so each country is updated exactly once.
$colIndex = @{ A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8 }

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $countryName = $ws.Cells.Item($r, 1).Value2
    if ($updates.ContainsKey($countryName)) {
        $cols = $updates[$countryName]
        foreach ($col in $cols.Keys) {
            $ws.Cells.Item($r, $colIndex[$col]).Value = $cols[$col]
        }
    }
}

# Re-sort the data table by "Casos totales" (column B) descending, like
# the source dashboard re-export does.
$dataRange = $ws.Range("A$($firstRow):H$($lastRow)")
$sortKey = $ws.Range("B$($firstRow):B$($lastRow)")
$dataRange.Sort($sortKey, 2)

# Bump the "updated at" footer timestamp.
$ws.Range("A1").Value = "Datos actualizados a 25 de Junio de 2020 a las 11:37"
